$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1510.6207
$ws.Range("J17").Value = 1651.1666
$ws.Range("L17").Value = 4953.4998
$ws.Range("N17").Value = -5289.4998
$ws.Range("H64").Value = 2829.375
$ws.Range("I64").Value = 2840
$ws.Range("J64").Value = 2824.5454
$ws.Range("K64").Value = 2840
$ws.Range("L64").Value = 2824.5454
$ws.Range("M64").Value = -2592
$ws.Range("N64").Value = -3320.5454
$ws.Range("H67").Value = 2829.375
$ws.Range("I67").Value = 2840
$ws.Range("J67").Value = 2824.5454
$ws.Range("K67").Value = 2840
$ws.Range("L67").Value = 2824.5454
$ws.Range("M67").Value = -1982
$ws.Range("N67").Value = -4540.5454
$ws.Range("H137").Value = 1231.7441
$ws.Range("I137").Value = 982.1667
$ws.Range("J137").Value = 1807.6923
$ws.Range("K137").Value = 2946.5001
$ws.Range("L137").Value = 5423.0769
$ws.Range("M137").Value = -396.5001000000002
$ws.Range("N137").Value = -10523.0769
$ws.Range("H138").Value = 2447.1475
$ws.Range("I138").Value = 1320.0555
$ws.Range("J138").Value = 3134.8645
$ws.Range("K138").Value = 3960.1665
$ws.Range("L138").Value = 9404.593500000001
$ws.Range("M138").Value = 1179.8335
$ws.Range("N138").Value = -19684.5935
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18479.8
$ws.Range("I32").Value = 18228.393
$ws.Range("K32").Value = 18228.393
$ws.Range("M32").Value = -17941.393
$ws.Range("H61").Value = 1679.2821
$ws.Range("J61").Value = 2277.0833
$ws.Range("L61").Value = 2277.0833
$ws.Range("N61").Value = -2701.0833
$ws.Range("H74").Value = 812.4912
$ws.Range("I74").Value = 812.94543
$ws.Range("J74").Value = 800
$ws.Range("K74").Value = 812.94543
$ws.Range("L74").Value = 800
$ws.Range("M74").Value = 61.05457000000001
$ws.Range("N74").Value = -2548
$ws.Range("H77").Value = 812.4912
$ws.Range("I77").Value = 812.94543
$ws.Range("J77").Value = 800
$ws.Range("K77").Value = 4064.72715
$ws.Range("L77").Value = 4000
$ws.Range("M77").Value = 303.2728500000003
$ws.Range("N77").Value = -12736
$ws.Range("H102").Value = 1312.7273
$ws.Range("I102").Value = 1244
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 1244
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = 378
$ws.Range("N102").Value = -5244
$ws.Range("H122").Value = 2071.963
$ws.Range("I122").Value = 1862.9048
$ws.Range("J122").Value = 2803.6667
$ws.Range("K122").Value = 5588.7144
$ws.Range("L122").Value = 8411.000100000001
$ws.Range("M122").Value = -3138.7144
$ws.Range("N122").Value = -13311.0001
$ws.Range("H132").Value = 7700.027
$ws.Range("I132").Value = 8743.931
$ws.Range("K132").Value = 26231.793
$ws.Range("M132").Value = -23701.793
$ws.Range("H136").Value = 1679.2821
$ws.Range("J136").Value = 2277.0833
$ws.Range("L136").Value = 6831.249899999999
$ws.Range("N136").Value = -11931.2499
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2275
$ws.Range("I105").Value = 2400
$ws.Range("K105").Value = 2400
$ws.Range("M105").Value = -653
$ws.Range("H134").Value = 4537.488
$ws.Range("I134").Value = 4898.3438
$ws.Range("J134").Value = 3254.4443
$ws.Range("K134").Value = 14695.0314
$ws.Range("L134").Value = 9763.332900000001
$ws.Range("M134").Value = -12160.0314
$ws.Range("N134").Value = -14833.3329
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5301.6
$ws.Range("I31").Value = 4302.857
$ws.Range("J31").Value = 6572.727
$ws.Range("K31").Value = 4302.857
$ws.Range("L31").Value = 6572.727
$ws.Range("M31").Value = -4007.857
$ws.Range("N31").Value = -7162.727
$ws.Range("H34").Value = 5301.6
$ws.Range("I34").Value = 4302.857
$ws.Range("J34").Value = 6572.727
$ws.Range("K34").Value = 4302.857
$ws.Range("L34").Value = 6572.727
$ws.Range("M34").Value = -4100.857
$ws.Range("N34").Value = -6976.727
$ws.Range("H58").Value = 1470.1666
$ws.Range("I58").Value = 1544.963
$ws.Range("K58").Value = 1544.963
$ws.Range("M58").Value = -1341.963
$ws.Range("H99").Value = 2406.25
$ws.Range("I99").Value = 2414.2856
$ws.Range("J99").Value = 2400
$ws.Range("K99").Value = 2414.2856
$ws.Range("L99").Value = 2400
$ws.Range("M99").Value = -916.2856000000002
$ws.Range("N99").Value = -5396
$ws.Range("H126").Value = 2406.25
$ws.Range("I126").Value = 2414.2856
$ws.Range("J126").Value = 2400
$ws.Range("K126").Value = 7242.8568
$ws.Range("L126").Value = 7200
$ws.Range("M126").Value = -4772.8568
$ws.Range("N126").Value = -12140
$ws.Range("H134").Value = 5502.727
$ws.Range("I134").Value = 6148.421
$ws.Range("K134").Value = 18445.263
$ws.Range("M134").Value = -15910.263
$ws.Range("H136").Value = 1470.1666
$ws.Range("I136").Value = 1544.963
$ws.Range("K136").Value = 4634.889
$ws.Range("M136").Value = -2084.889
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 2429.5557
$ws.Range("I81").Value = 392.66666
$ws.Range("J81").Value = 3448
$ws.Range("K81").Value = 1177.99998
$ws.Range("L81").Value = 10344
$ws.Range("M81").Value = -54.99998000000005
$ws.Range("N81").Value = -12590
$ws.Range("H84").Value = 2429.5557
$ws.Range("I84").Value = 392.66666
$ws.Range("J84").Value = 3448
$ws.Range("K84").Value = 3533.99994
$ws.Range("L84").Value = 31032
$ws.Range("M84").Value = 2082.00006
$ws.Range("N84").Value = -42264
$ws.Range("H131").Value = 591.74
$ws.Range("I131").Value = 256.53845
$ws.Range("J131").Value = 806.0492
$ws.Range("K131").Value = 769.61535
$ws.Range("L131").Value = 2418.1476
$ws.Range("M131").Value = 4270.38465
$ws.Range("N131").Value = -12498.1476
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1858.7778
$ws.Range("I97").Value = 2043
$ws.Range("J97").Value = 1379.8
$ws.Range("K97").Value = 2043
$ws.Range("L97").Value = 1379.8
$ws.Range("M97").Value = -1547
$ws.Range("N97").Value = -2371.8
$ws.Range("H102").Value = 1775.0714
$ws.Range("I102").Value = 1266.3334
$ws.Range("J102").Value = 2690.8
$ws.Range("K102").Value = 1266.3334
$ws.Range("L102").Value = 2690.8
$ws.Range("M102").Value = 355.6666
$ws.Range("N102").Value = -5934.8
$ws.Range("H122").Value = 27779980
$ws.Range("I122").Value = 37038972
$ws.Range("K122").Value = 111116916
$ws.Range("M122").Value = -111114466
$ws.Range("H123").Value = 28175.584
$ws.Range("J123").Value = 28175.584
$ws.Range("L123").Value = 28175.584
$ws.Range("N123").Value = -33075.584
$ws.Range("H132").Value = 5301.161
$ws.Range("I132").Value = 5733.52
$ws.Range("J132").Value = 3499.6667
$ws.Range("K132").Value = 17200.56
$ws.Range("L132").Value = 10499.0001
$ws.Range("M132").Value = -14670.56
$ws.Range("N132").Value = -15559.0001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 433.9
$ws.Range("I22").Value = 919.6667
$ws.Range("K22").Value = 919.6667
$ws.Range("M22").Value = -624.6667
$ws.Range("H27").Value = 433.9
$ws.Range("I27").Value = 919.6667
$ws.Range("K27").Value = 919.6667
$ws.Range("M27").Value = -812.6667
$ws.Range("H46").Value = 1230
$ws.Range("I46").Value = 549.6667
$ws.Range("J46").Value = 1407.4783
$ws.Range("K46").Value = 549.6667
$ws.Range("L46").Value = 1407.4783
$ws.Range("M46").Value = -361.6667
$ws.Range("N46").Value = -1783.4783
$ws.Range("H61").Value = 20834678
$ws.Range("I61").Value = 1288
$ws.Range("J61").Value = 47620464
$ws.Range("K61").Value = 1288
$ws.Range("L61").Value = 47620464
$ws.Range("M61").Value = -1086
$ws.Range("N61").Value = -47620868
$ws.Range("H68").Value = 1304.5
$ws.Range("I68").Value = 1056
$ws.Range("J68").Value = 1884.3334
$ws.Range("K68").Value = 1056
$ws.Range("L68").Value = 1884.3334
$ws.Range("M68").Value = -307
$ws.Range("N68").Value = -3382.3334
$ws.Range("H71").Value = 1304.5
$ws.Range("I71").Value = 1056
$ws.Range("J71").Value = 1884.3334
$ws.Range("K71").Value = 5280
$ws.Range("L71").Value = 9421.666999999999
$ws.Range("M71").Value = -1536
$ws.Range("N71").Value = -16909.667
$ws.Range("H82").Value = 1120.9231
$ws.Range("I82").Value = 1078.9166
$ws.Range("K82").Value = 1078.9166
$ws.Range("M82").Value = -717.9166
$ws.Range("H85").Value = 1120.9231
$ws.Range("I85").Value = 1078.9166
$ws.Range("K85").Value = 1078.9166
$ws.Range("M85").Value = 169.0834
$ws.Range("H100").Value = 66667170
$ws.Range("I100").Value = 66667170
$ws.Range("K100").Value = 66667170
$ws.Range("M100").Value = -66666629
$ws.Range("H113").Value = 20834678
$ws.Range("I113").Value = 1288
$ws.Range("J113").Value = 47620464
$ws.Range("K113").Value = 1288
$ws.Range("L113").Value = 47620464
$ws.Range("M113").Value = 882
$ws.Range("N113").Value = -47624804
$ws.Range("H132").Value = 8388.9375
$ws.Range("I132").Value = 11071.857
$ws.Range("J132").Value = 3267
$ws.Range("K132").Value = 33215.571
$ws.Range("L132").Value = 9801
$ws.Range("M132").Value = -30685.571
$ws.Range("N132").Value = -14861
$ws.Range("H133").Value = 24819.385
$ws.Range("J133").Value = 24819.385
$ws.Range("L133").Value = 24819.385
$ws.Range("N133").Value = -29879.385
$ws.Range("H136").Value = 5294.375
$ws.Range("I136").Value = 5857.2
$ws.Range("K136").Value = 17571.6
$ws.Range("M136").Value = -15021.6
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 12000
$ws.Range("I14").Value = 12000
$ws.Range("K14").Value = 12000
$ws.Range("M14").Value = -11832
$ws.Range("H132").Value = 1954.5682
$ws.Range("I132").Value = 1835.2354
$ws.Range("J132").Value = 2360.3
$ws.Range("K132").Value = 5505.706200000001
$ws.Range("L132").Value = 7080.900000000001
$ws.Range("M132").Value = -2975.706200000001
$ws.Range("N132").Value = -12140.9
